$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last timestamp column (GM) into a brand new column so that the
# new column inherits the same per-row cell typing (numeric price vs. blank)
# as the existing GM column, then insert it right before the old GN column.
# This shifts the old GN (nom) and GO (url_produit) columns one place to the
# right, becoming GO and GP respectively - matching the diff.
$ws.Columns("GM").Copy() | Out-Null
$ws.Columns("GN").Insert()

# The freshly inserted column copied GM's timestamp header; replace it with
# the new scrape timestamp for this run.
$ws.Range("GN1").Value = "2026-02-06 00:58:30"
